# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (rows 16-52, column E) is re-sorted from
# descending (2003 .. 1703) to ascending (1703 .. 2003) order, and the
# "Valor Mora" (F) / "Salario Basico" (G) figures for each period are
# refreshed to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New ascending period order (previously descending) for rows 16..52
$periods = @(
    "1703", "1704", "1705", "1706", "1707", "1708", "1709", "1710", "1711", "1712",
    "1801", "1802", "1803", "1804", "1805", "1806", "1807", "1808", "1809", "1810",
    "1811", "1812", "1901", "1902", "1903", "1904", "1905", "1906", "1907", "1908",
    "1909", "1910", "1911", "1912", "2001", "2002", "2003"
)

$startRow = 16

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i

    # Column E: period label
    $ws.Range("E$row").Value = $periods[$i]

    # Column F: Valor Mora - first 18 periods (1703..1808) use 27578,
    # the remaining 19 periods (1809..2003) use 31249
    if ($i -lt 18) {
        $ws.Range("F$row").Value = 27578
    } else {
        $ws.Range("F$row").Value = 31249
    }

    # Column G: Salario Basico - updated uniformly to 781242
    $ws.Range("G$row").Value = 781242
}
